$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.581554651260376
$ws.Range("B1").Value = 1.954456090927124
$ws.Range("C1").Value = 2.18241810798645
$ws.Range("D1").Value = 2.413662672042847
$ws.Range("E1").Value = 3.079207897186279
